$d = $word.ActiveDocument

# 1. Replace the placeholder name with the student's name, keeping the
#    surrounding bracket characters and yellow-highlight formatting intact,
#    but splitting the bracketed text into three runs ("[", name, "]") to
#    match the authored edit.
$rng = $d.Content
$find = $rng.Find
$find.ClearFormatting()
$find.Text = "[Insert your name here]"
$found = $find.Execute()

$fullStart = $rng.Start
$fullEnd = $rng.End

$nameRange = $d.Range($fullStart + 1, $fullEnd - 1)
$newName = "Tiago Neves Sousa"
$nameRange.Text = $newName

$nameRange2 = $d.Range($fullStart + 1, $fullStart + 1 + $newName.Length)
$nameRange2.Bold = 1
$nameRange2.Bold = 0

# 2. Fill in the first data row of the "Professional Development Plan"
#    table (Rating, Week #, Description columns) for the first outcome
#    row ("Become more efficient at applying your innate curiosity and
#    creativity.").
$table = $d.Tables(3)
$table.Cell(2, 2).Range.Text = "3"
$table.Cell(2, 3).Range.Text = "2"
$table.Cell(2, 4).Range.Text = "I" + [char]8217 + "m trying to organize my life to meet with teammates and have done some tasks solution by myself. Also joined the team Trello and Github to accomplish the material requests."
